$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 44 (pushes rows 44-118 down to 45-119), inheriting the
# formatting of row 43 (the row immediately above) for the new row.
$ws.Rows(44).Insert()

# Populate the new row with the Juba/Khartoum (hj-hs) entry.
$ws.Range("A44").Value = "hj-hs"
$ws.Range("B44").Value = "Juba, Khartoum"
$ws.Range("C44").Value = "VATGlasses"
$ws.Range("D44").Value = "vACC Rejected Offer To Participate"

# Merge D44:E44 like the surrounding "Manager" cells.
$ws.Range("D44:E44").Merge()

# Restore selection/scroll state similar to the post-edit workbook.
$ws.Range("D44:E44").Select()
